$d = $word.ActiveDocument

# The date cell currently contains two runs: "03October" and "2022".
# Together they read "03October2022" with no space between them. Replace
# that whole string with the single new date "1 November 2023".

$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute("03October2022", $true, $false, $false, $false, $false, `
               $true, 1, $false, "1 November 2023", 2)
